# "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
#
# Rebuilds the worker/period "Estado de Cuenta" detail rows (B16:J35) so the
# records are grouped by period (2408, 2409, 2410, 2411, 2412, 2501) instead
# of by worker, and adds the new period rows (2408-2412) that didn't exist
# before. Column F (Valor Mora) is 52000 for every period except 2501, which
# keeps 50266.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tipo Doc (col B) is "CC" for every detail row and Salario Basico (col G)
# is always 1300000 - neither changes, only C:F below vary.
$workers = @(
    @{ Doc = "9158235";  Name = "FERNANDO JOSE JIMENEZ BARRETO" },
    @{ Doc = "32907838"; Name = "BANIDIS DEL CARMEN BANQUET BLANCO" },
    @{ Doc = "45579779"; Name = "GLADYS DEL SOCORRO LEGUIA ROBLES" },
    @{ Doc = "73121215"; Name = "YIDIO MORALES TORRES" }
)

# Period -> which workers (by index into $workers) have a record that period,
# in the row order they appear in the sheet.
$periods = @(
    @{ Period = "2408"; Members = @(0,1,2,3); Valor = 52000 },
    @{ Period = "2409"; Members = @(0,1,2,3); Valor = 52000 },
    @{ Period = "2410"; Members = @(1,2);     Valor = 52000 },
    @{ Period = "2411"; Members = @(1,2);     Valor = 52000 },
    @{ Period = "2412"; Members = @(0,1,2,3); Valor = 52000 },
    @{ Period = "2501"; Members = @(0,1,2,3); Valor = 50266 }
)

$row = 16
foreach ($p in $periods) {
    foreach ($idx in $p.Members) {
        $w = $workers[$idx]
        $ws.Range("C$row").Value = $w.Doc
        $ws.Range("D$row").Value = $w.Name
        $ws.Range("E$row").Value = $p.Period
        $ws.Range("F$row").Value = $p.Valor
        $ws.Range("G$row").Value = 1300000
        $row++
    }
}
